$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "PROFILE" column header in E1, matching the header style of B1:D1
$ws.Range("E1").Value = "PROFILE"
$ws.Range("E1").Font.Size = $ws.Range("D1").Font.Size
$ws.Range("E1").HorizontalAlignment = $ws.Range("D1").HorizontalAlignment

# Fill in profile values for each user row, matching the style of the data cells in C:D
$ws.Range("E2").Value = "admin"
$ws.Range("E2").Font.Size = $ws.Range("C2").Font.Size
$ws.Range("E2").HorizontalAlignment = $ws.Range("C2").HorizontalAlignment

$ws.Range("E3").Value = "editor"
$ws.Range("E3").Font.Size = $ws.Range("C3").Font.Size
$ws.Range("E3").HorizontalAlignment = $ws.Range("C3").HorizontalAlignment

$ws.Range("E4").Value = "editor"
$ws.Range("E4").Font.Size = $ws.Range("C4").Font.Size
$ws.Range("E4").HorizontalAlignment = $ws.Range("C4").HorizontalAlignment

# Match new column width for column E (target stored width ~22.42578125 chars)
$ws.Columns.Item(5).ColumnWidth = 21.67

# Update the selected cell to reflect where the user ended up after editing
$ws.Range("E5").Select()
